$d = $word.ActiveDocument

$label = "MANUSCRIPT NUMBER"

# Replace "FULL ARTICLE DOI" with "MANUSCRIPT NUMBER"
$d.Content.Find.Execute("FULL ARTICLE DOI", $true, $false, $false, $false, $false,
                         $true, 1, $false, $label, 2)

# Replace the DOI URL with just the manuscript number "99999"
$d.Content.Find.Execute("https://doi.org/10.7554/eLife.99999", $true, $false, $false, $false, $false,
                         $true, 1, $false, "99999", 2)

# The paragraph now reads "MANUSCRIPT NUMBER" + line break + "99999", all
# originally inside a single run for the label+break. Split the line break
# off into its own run (matching the target OOXML) by nudging a character
# formatting property on just that break position.
$para = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $candidate = $d.Paragraphs($i)
    if ($candidate.Range.Text.StartsWith($label)) {
        $para = $candidate
    }
}

$brPos = $para.Range.Start + $label.Length
$brRange = $d.Range($brPos, $brPos + 1)
$brRange.Font.Bold = $false
$brRange.Font.Bold = $true
